$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "('Bat', ['Token Creature — Bat', 'Flying', '1/1'])",
    "('Cat', ['Token Creature — Cat', '2/2'])",
    "('Cat Dragon', ['Token Creature — Cat Dragon', 'Flying', '3/3'])",
    "('Cat Warrior', ['Token Creature — Cat Warrior', 'Forestwalk', '2/2'])",
    "('Dragon', ['Token Creature — Dragon', 'Flying', '6/6'])",
    "('Eldrazi Spawn', ['Token Creature — Eldrazi Spawn', 'Sacrifice this creature: Add {C}.', '0/1'])",
    "('Gold', ['Token Artifact — Gold', 'Sacrifice this artifact: Add one mana of any color.'])",
    "('Rat', ['Token Creature — Rat', 'Deathtouch', '1/1'])",
    "('Vampire', ['Token Creature — Vampire', '1/1'])",
    "('Zombie', ['Token Creature — Zombie', '2/2'])"
)

# Clear old rows beyond the new range (rows 2 through 41)
$ws.Range("A2:A41").ClearContents()

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
